# Add a new order line (row 19) to the UNFI Triphammer order sheet.
# Values must be written as text (matching the existing inlineStr cells),
# not auto-coerced into numbers/currency by Excel, so we build each cell
# via a text formula ("=""literal""") and then convert the formulas to
# static values with a copy / paste-special(values) round trip.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Formula = '="0778936"'
$ws.Range("B19").Formula = '="Seeds - Pumpkin Seeds Raw"'
$ws.Range("C19").Formula = '="2"'
$ws.Range("D19").Formula = '="$32.53"'
$ws.Range("E19").Formula = '="$65.06"'

$ws.Range("A19:E19").Copy()
$ws.Range("A19:E19").PasteSpecial(-4163)
